$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.933.91'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.506.09'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.35'
$ws.Range('E5').Value = '  -1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '195.96'
$ws.Range('E6').Value = '  +6.44%  '
$ws.Range('E7').Value = '  +1.77%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('E10').Value = '  +2.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.09'
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('E12').Value = '  -2.12%  '
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.059.41'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '605.97'
$ws.Range('E15').Value = '  +3.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '70.040.54'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.60'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.503.05'
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('E20').Value = '  +0.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.993'
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.07'
$ws.Range('E22').Value = '  +3.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '104.33'
$ws.Range('E23').Value = '  +8.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.11'
$ws.Range('E24').Value = '  +5.72%  '
$ws.Range('E25').Value = '  -2.09%  '
$ws.Range('E26').Value = '  +3.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.71'
$ws.Range('E28').Value = '  +2.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.57'
$ws.Range('E29').Value = '  +5.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.59'
$ws.Range('E30').Value = '  +27.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.11'
$ws.Range('E31').Value = '  +1.95%  '
$ws.Range('E32').Value = '  +4.31%  '
$ws.Range('E33').Value = '  +1.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.22'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0823'
$ws.Range('E35').Value = '  +6.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.725.40'
$ws.Range('E36').Value = '  +5.55%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.08'
$ws.Range('E37').Value = '  -5.67%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.393'
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.83'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.58'
$ws.Range('E41').Value = '  +2.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '501.04'
$ws.Range('E42').Value = '  -5.54%  '
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0458'
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.33'
$ws.Range('E45').Value = '  -0.71%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.140'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.83'
$ws.Range('E47').Value = '  -3.49%  '
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('E49').Value = '  -4.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000246'
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '131.41'
$ws.Range('E51').Value = '  -2.99%  '
